$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The whole data table (rows 2-13, years 2009..2020) shifts up by one year:
# the old 2009 row is removed and a new 2021 row is appended at the bottom.
# Deleting row 2 (2009) shifts every following row up by one, which matches
# rows 2..12 of the target (2010..2020) exactly, including carrying along
# the untouched/blank cells and existing formatting.
$ws.Rows.Item(2).Delete()

# Now append the new last row (13) for 2021.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = 71.9406
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = 221.1286
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = 816.2542
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = 169.5455
$ws.Range("J13").Value = 113.1454
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = 0.0002
$ws.Range("O13").Value = 185.086
$ws.Range("P13").Value = ""
$ws.Range("Q13").Value = ""
$ws.Range("R13").Value = ""
$ws.Range("S13").Value = ""

# Match the header-cell formatting used by A2:A12 (bold, thin border all
# around, centered / top-aligned) on the new A13 label cell.
$ws.Range("A13").Font.Bold = $true
$ws.Range("A13").HorizontalAlignment = -4108
$ws.Range("A13").VerticalAlignment = -4160
$ws.Range("A13").Borders.LineStyle = 1
